$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3333333333333333
$ws.Range("C2").Value = 0.3
$ws.Range("D2").Value = 0.3157894736842105

# Row 3
$ws.Range("B3").Value = 0.5333333333333333
$ws.Range("C3").Value = 0.5714285714285714
$ws.Range("D3").Value = 0.5517241379310344

# Row 4
$ws.Range("B4").Value = 0.4583333333333333
$ws.Range("C4").Value = 0.4583333333333333
$ws.Range("D4").Value = 0.4583333333333333
$ws.Range("E4").Value = 0.4583333333333333

# Row 5
$ws.Range("B5").Value = 0.4333333333333333
$ws.Range("C5").Value = 0.4357142857142857
$ws.Range("D5").Value = 0.4337568058076224

# Row 6
$ws.Range("B6").Value = 0.45
$ws.Range("C6").Value = 0.4583333333333333
$ws.Range("D6").Value = 0.4534180278281911

# Row 7
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = 0.5

# Row 8
$ws.Range("B8").Value = 0.6428571428571429
$ws.Range("C8").Value = 0.6428571428571429
$ws.Range("D8").Value = 0.6428571428571429

# Row 9
$ws.Range("B9").Value = 0.5833333333333334
$ws.Range("C9").Value = 0.5833333333333334
$ws.Range("D9").Value = 0.5833333333333334
$ws.Range("E9").Value = 0.5833333333333334

# Row 10
$ws.Range("B10").Value = 0.5714285714285714
$ws.Range("C10").Value = 0.5714285714285714
$ws.Range("D10").Value = 0.5714285714285714

# Row 11
$ws.Range("B11").Value = 0.5833333333333334
$ws.Range("C11").Value = 0.5833333333333334
$ws.Range("D11").Value = 0.5833333333333334

# Row 22
$ws.Range("B22").Value = 0.6
$ws.Range("C22").Value = 0.6
$ws.Range("D22").Value = 0.6

# Row 23
$ws.Range("B23").Value = 0.7142857142857143
$ws.Range("C23").Value = 0.7142857142857143
$ws.Range("D23").Value = 0.7142857142857143

# Row 24
$ws.Range("B24").Value = 0.6666666666666666
$ws.Range("C24").Value = 0.6666666666666666
$ws.Range("D24").Value = 0.6666666666666666
$ws.Range("E24").Value = 0.6666666666666666

# Row 25
$ws.Range("B25").Value = 0.6571428571428571
$ws.Range("C25").Value = 0.6571428571428571
$ws.Range("D25").Value = 0.6571428571428571

# Row 26
$ws.Range("B26").Value = 0.6666666666666666
$ws.Range("C26").Value = 0.6666666666666666
$ws.Range("D26").Value = 0.6666666666666666
